$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 360.2
$ws.Range("I15").Value = 360.2
$ws.Range("K15").Value = 1080.6
$ws.Range("M15").Value = -911.5999999999999
# Row 125
$ws.Range("H125").Value = 62500984
$ws.Range("I125").Value = 90909790
$ws.Range("J125").Value = 1591.8
$ws.Range("K125").Value = 818188110
$ws.Range("L125").Value = 14326.2
$ws.Range("M125").Value = -818185650
$ws.Range("N125").Value = -19246.2
# Row 137
$ws.Range("H137").Value = 798.775
$ws.Range("I137").Value = 698.8095
$ws.Range("K137").Value = 2096.4285
$ws.Range("M137").Value = 453.5715
# Row 138
$ws.Range("H138").Value = 2598.2354
$ws.Range("I138").Value = 1166.2391
$ws.Range("J138").Value = 5592.409
$ws.Range("K138").Value = 3498.7173
$ws.Range("L138").Value = 16777.227
$ws.Range("M138").Value = 1641.2827
$ws.Range("N138").Value = -27057.227

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 384117.25
$ws.Range("I32").Value = 2664
$ws.Range("J32").Value = 4089663.2
$ws.Range("K32").Value = 2664
$ws.Range("L32").Value = 4089663.2
$ws.Range("M32").Value = -2377
$ws.Range("N32").Value = -4090237.2
# Row 37
$ws.Range("H37").Value = 250006140
$ws.Range("I37").Value = 1000000000
$ws.Range("J37").Value = 8201
$ws.Range("K37").Value = 1000000000
$ws.Range("L37").Value = 8201
$ws.Range("M37").Value = -999999727
$ws.Range("N37").Value = -8747
# Row 74
$ws.Range("H74").Value = 213620.28
$ws.Range("I74").Value = 228097.8
$ws.Range("J74").Value = 1283.3334
$ws.Range("K74").Value = 228097.8
$ws.Range("L74").Value = 1283.3334
$ws.Range("M74").Value = -227223.8
$ws.Range("N74").Value = -3031.3334
# Row 77
$ws.Range("H77").Value = 213620.28
$ws.Range("I77").Value = 228097.8
$ws.Range("J77").Value = 1283.3334
$ws.Range("K77").Value = 1140489
$ws.Range("L77").Value = 6416.666999999999
$ws.Range("M77").Value = -1136121
$ws.Range("N77").Value = -15152.667
# Row 102
$ws.Range("H102").Value = 3863.389
$ws.Range("I102").Value = 2902
$ws.Range("K102").Value = 2902
$ws.Range("M102").Value = -1280
# Row 132
$ws.Range("H132").Value = 1057.6666
$ws.Range("I132").Value = 781.1627999999999
$ws.Range("J132").Value = 3435.6
$ws.Range("K132").Value = 2343.4884
$ws.Range("L132").Value = 10306.8
$ws.Range("M132").Value = 186.5116000000003
$ws.Range("N132").Value = -15366.8

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 82
$ws.Range("H82").Value = 18616.6
$ws.Range("I82").Value = 8000
$ws.Range("J82").Value = 25694.334
$ws.Range("K82").Value = 8000
$ws.Range("L82").Value = 25694.334
$ws.Range("M82").Value = -7617
$ws.Range("N82").Value = -26460.334
# Row 85
$ws.Range("H85").Value = 18616.6
$ws.Range("I85").Value = 8000
$ws.Range("J85").Value = 25694.334
$ws.Range("K85").Value = 8000
$ws.Range("L85").Value = 25694.334
$ws.Range("M85").Value = -6674
$ws.Range("N85").Value = -28346.334
# Row 134
$ws.Range("H134").Value = 821.2353000000001
$ws.Range("I134").Value = 613.2727
$ws.Range("J134").Value = 2128.4285
$ws.Range("K134").Value = 1839.8181
$ws.Range("L134").Value = 6385.2855
$ws.Range("M134").Value = 695.1819
$ws.Range("N134").Value = -11455.2855

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 6851035
$ws.Range("I31").Value = 8065825
$ws.Range("K31").Value = 8065825
$ws.Range("M31").Value = -8065530
# Row 34
$ws.Range("H34").Value = 6851035
$ws.Range("I34").Value = 8065825
$ws.Range("K34").Value = 8065825
$ws.Range("M34").Value = -8065623
# Row 50
$ws.Range("H50").Value = 9401
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 9401
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 9401
$ws.Range("M50").Value = $null
$ws.Range("N50").Value = -10651
# Row 51
$ws.Range("H51").Value = 9500.75
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 9500.75
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 9500.75
$ws.Range("M51").Value = $null
$ws.Range("N51").Value = -10972.75
# Row 58
$ws.Range("H58").Value = 1087.1842
$ws.Range("I58").Value = 987.8095
$ws.Range("J58").Value = 1209.9412
$ws.Range("K58").Value = 987.8095
$ws.Range("L58").Value = 1209.9412
$ws.Range("M58").Value = -784.8095
$ws.Range("N58").Value = -1615.9412
# Row 60
$ws.Range("H60").Value = 8735
$ws.Range("J60").Value = 8735
$ws.Range("L60").Value = 8735
$ws.Range("N60").Value = -9757
# Row 61
$ws.Range("H61").Value = 9500.75
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 9500.75
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 9500.75
$ws.Range("M61").Value = $null
$ws.Range("N61").Value = -10196.75
# Row 80
$ws.Range("H80").Value = 27800
$ws.Range("J80").Value = 27800
$ws.Range("L80").Value = 27800
$ws.Range("N80").Value = -30046
# Row 83
$ws.Range("H83").Value = 27800
$ws.Range("J83").Value = 27800
$ws.Range("L83").Value = 83400
$ws.Range("N83").Value = -94632
# Row 136
$ws.Range("H136").Value = 1087.1842
$ws.Range("I136").Value = 987.8095
$ws.Range("J136").Value = 1209.9412
$ws.Range("K136").Value = 2963.4285
$ws.Range("L136").Value = 3629.8236
$ws.Range("M136").Value = -413.4285
$ws.Range("N136").Value = -8729.8236

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1045.75
$ws.Range("I5").Value = 834.3333
$ws.Range("J5").Value = 1680
$ws.Range("K5").Value = 2502.9999
$ws.Range("L5").Value = 5040
$ws.Range("M5").Value = -2390.9999
$ws.Range("N5").Value = -5264
# Row 113
$ws.Range("H113").Value = 984.4167
$ws.Range("I113").Value = 761.5
$ws.Range("J113").Value = 1018.71155
$ws.Range("K113").Value = 2284.5
$ws.Range("L113").Value = 3056.13465
$ws.Range("M113").Value = -114.5
$ws.Range("N113").Value = -7396.13465
# Row 131
$ws.Range("H131").Value = 7693165.5
$ws.Range("I131").Value = 771.3570999999999
$ws.Range("J131").Value = 9804803
$ws.Range("K131").Value = 2314.0713
$ws.Range("L131").Value = 29414409
$ws.Range("M131").Value = 2725.9287
$ws.Range("N131").Value = -29424489
# Row 135
$ws.Range("H135").Value = 1045.75
$ws.Range("I135").Value = 834.3333
$ws.Range("J135").Value = 1680
$ws.Range("K135").Value = 7508.9997
$ws.Range("L135").Value = 15120
$ws.Range("M135").Value = -4973.9997
$ws.Range("N135").Value = -20190

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 14429947
$ws.Range("I70").Value = 26792010
$ws.Range("J70").Value = 7541.222
$ws.Range("K70").Value = 26792010
$ws.Range("L70").Value = 7541.222
$ws.Range("M70").Value = -26791740
$ws.Range("N70").Value = -8081.222
# Row 73
$ws.Range("H73").Value = 14429947
$ws.Range("I73").Value = 26792010
$ws.Range("J73").Value = 7541.222
$ws.Range("K73").Value = 26792010
$ws.Range("L73").Value = 7541.222
$ws.Range("M73").Value = -26791074
$ws.Range("N73").Value = -9413.222
# Row 97
$ws.Range("H97").Value = 1216.3043
$ws.Range("I97").Value = 720
$ws.Range("J97").Value = 1988.3334
$ws.Range("K97").Value = 720
$ws.Range("L97").Value = 1988.3334
$ws.Range("M97").Value = -224
$ws.Range("N97").Value = -2980.3334
# Row 102
$ws.Range("H102").Value = 1353.7727
$ws.Range("I102").Value = 1268.5278
$ws.Range("J102").Value = 1737.375
$ws.Range("K102").Value = 1268.5278
$ws.Range("L102").Value = 1737.375
$ws.Range("M102").Value = 353.4721999999999
$ws.Range("N102").Value = -4981.375
# Row 132
$ws.Range("H132").Value = 21304.549
$ws.Range("I132").Value = 26699.15
$ws.Range("J132").Value = 1687.8182
$ws.Range("K132").Value = 80097.45000000001
$ws.Range("L132").Value = 5063.4546
$ws.Range("M132").Value = -77567.45000000001
$ws.Range("N132").Value = -10123.4546

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 1902.8235
$ws.Range("I68").Value = 1800.3
$ws.Range("J68").Value = 2049.2856
$ws.Range("K68").Value = 1800.3
$ws.Range("L68").Value = 2049.2856
$ws.Range("M68").Value = -1051.3
$ws.Range("N68").Value = -3547.2856
# Row 71
$ws.Range("H71").Value = 1902.8235
$ws.Range("I71").Value = 1800.3
$ws.Range("J71").Value = 2049.2856
$ws.Range("K71").Value = 9001.5
$ws.Range("L71").Value = 10246.428
$ws.Range("M71").Value = -5257.5
$ws.Range("N71").Value = -17734.428
# Row 100
$ws.Range("H100").Value = 2257.3872
$ws.Range("I100").Value = 1767.5264
$ws.Range("J100").Value = 3033
$ws.Range("K100").Value = 1767.5264
$ws.Range("L100").Value = 3033
$ws.Range("M100").Value = -1226.5264
$ws.Range("N100").Value = -4115
# Row 132
$ws.Range("H132").Value = 2610.7593
$ws.Range("I132").Value = 2829.1707
$ws.Range("J132").Value = 1921.9231
$ws.Range("K132").Value = 8487.5121
$ws.Range("L132").Value = 5765.7693
$ws.Range("M132").Value = -5957.5121
$ws.Range("N132").Value = -10825.7693

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 123
$ws.Range("H123").Value = 32663.334
$ws.Range("J123").Value = 32663.334
$ws.Range("L123").Value = 32663.334
$ws.Range("N123").Value = -42463.334
# Row 132
$ws.Range("H132").Value = 29765648
$ws.Range("I132").Value = 48078424
$ws.Range("J132").Value = 7385
$ws.Range("K132").Value = 144235272
$ws.Range("L132").Value = 22155
$ws.Range("M132").Value = -144232742
$ws.Range("N132").Value = -27215
